# UserController erweitert und ChangeUserData.cshtml erstellt
#
# Tasks sheet: mark the "ChangeUserData" sub-tasks as done (with completion
# dates), add a new "Navigations-Menü für Zugriff auf das Profil erweitern"
# task, and mark the "IUserRepository ... (ChangePassword)" task as still
# in progress ("b").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# Row 26: "IUserRepository und UserRepositoryDB erweitern (ChangeUserData und GetUser)"
# used to be marked "b" (in Bearbeitung) - now it's done, with a completion date.
$ws.Range("C26").ClearContents()
$ws.Range("C26").Value = "done"
$src = $ws.Range("D4")
$src.Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("D26").Value2 = 43530

# Row 27: "ChangeUserData Methode im UserController erstellen" -> done
$ws.Range("C27").Value = "done"
$ws.Range("D4").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value2 = 43531

# Row 28: "ChangeUserData View erstellen" -> done
$ws.Range("C28").Value = "done"
$ws.Range("D4").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value2 = 43531

# Insert a new row 29 for the new task "Navigations-Menü für Zugriff auf das
# Profil erweitern", also marked done.
$ws.Rows("29").Insert()
$ws.Range("B29").Value = "Navigations-Menü für Zugriff auf das Profil erweitern"
$ws.Range("C29").Value = "done"
$ws.Range("D4").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").Value2 = 43531

# The "IUserRepository ... (ChangePassword)" task (now shifted to row 30) is
# still in progress.
$ws.Range("C30").Value = "b"

$excel.CutCopyMode = $false

# Keep the active selection in sync with the newly added row (matches the
# author ending up on the new ChangePassword "b" cell after the edit).
$ws.Range("C30").Select()
